$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B39 was mistakenly stored as a text "2" - fix it to be a real number 2
$ws.Range("B39").Value = 2

# Add new row 40 with the new annotation data.
# B40 holds "4" as text (matching the surrounding inline-string cells), so
# force the Text number format before assignment to stop auto-numeric coercion.
$ws.Range("A40").Value = "Ying Tang"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "4"
$ws.Range("C40").Value = "motivated well,departs from prior work"
$ws.Range("D40").Value = "FBK"
$ws.Range("E40").Value = "THE"
$ws.Range("F40").Value = "d0296b92-10f5-497e-8726-aae675ac805b"
$ws.Range("G40").Value = "rJl3yM-Ab_annotated.xlsx"
$ws.Range("H40").Value = "The new method is motivated well and departs from prior work."
